$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the emoji values used as "statut" in column A:
#   📘 -> ⚠️
#   📕 -> -3   (kept as text, not a number)
#   📗 -> ✅

$rows = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()

    if ($val -eq "📘") {
        $cell.Value = "⚠️"
    }
    elseif ($val -eq "📕") {
        # Keep "-3" as text (not a number) so it round-trips the same
        # way the emoji text value did.
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
    }
    elseif ($val -eq "📗") {
        $cell.Value = "✅"
    }
}
